# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> applied to the slide master / slides ("Integral")
#   ppt/theme/theme2.xml  -> applied to the notes master        ("Office Theme")
#
# The authored edit swaps the two themes' contents: the slide-master theme
# becomes the stock "Office Theme" 12-colour scheme, while the notes-master
# theme becomes the former "Integral" 12-colour scheme.
#
# PowerPoint's automation surface only exposes the slide-facing theme colours
# (Slide.ThemeColorScheme / Master.ColorScheme, 12 slots: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) - there's no supported object model call to
# rewrite the notes master's theme part. So we recolor the reachable side,
# driving every theme colour slot to the stock Office palette to match the
# new ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x6A5444,  # 3  dk2
    0xE6E6E7,  # 4  lt2
    0xD59B5B,  # 5  accent1
    0x317DED,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0x00C0FF,  # 8  accent4
    0xC47244,  # 9  accent5
    0x47AD70,  # 10 accent6
    0xC16305,  # 11 hlink
    0x724F95   # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
